# EM_2021.xlsx update: EURO 2020 was postponed by exactly 364 days (one year,
# same weekday) to June 2021, and by the time of this edit the four playoff
# winners (paths A-D) were known, so the "Play-off X" placeholders used in
# the fixture list get replaced by the real qualified teams. The "Countries"
# reference sheet also drops the placeholder code "FYR_MACEDONIA" in favour
# of the plain "MACEDONIA" label and gets re-sorted so it lands in the right
# alphabetical spot.

$wb = $excel.ActiveWorkbook

# --- "Countries" sheet -----------------------------------------------------
# Rename the FYR_MACEDONIA entry to MACEDONIA, then re-sort the A1:A210
# list alphabetically so the renamed entry moves next to MACAU/MADAGASCAR.
$wsCountries = $wb.Worksheets.Item("Countries")
$wsCountries.Range("A71").Value = "MACEDONIA"

$wsCountries.Sort.SortFields.Clear()
$wsCountries.Sort.SortFields.Add($wsCountries.Range("A1:A210"))
$wsCountries.Sort.SetRange($wsCountries.Range("A1:A210"))
$wsCountries.Sort.Header = -4142
$wsCountries.Sort.Apply()

$wsCountries.Range("C12").Select()

# --- "Matches" sheet --------------------------------------------------------
# Shift every kickoff date by 364 days and swap the playoff placeholders for
# the real teams that qualified (Path A = Hungary, Path B = Slovakia,
# Path C = Scotland, Path D = Macedonia).
$ws = $wb.Worksheets.Item("Matches")

$ws.Range("D2").Value = 44358.875
$ws.Range("D3").Value = 44359.625
$ws.Range("D4").Value = 44363.875
$ws.Range("D5").Value = 44363.75
$ws.Range("D6").Value = 44367.75
$ws.Range("D7").Value = 44367.75

$ws.Range("D8").Value = 44359.875
$ws.Range("D9").Value = 44359.75
$ws.Range("D10").Value = 44363.625
$ws.Range("D11").Value = 44364.75
$ws.Range("D12").Value = 44368.875
$ws.Range("D13").Value = 44368.875

$ws.Range("D14").Value = 44360.875
$ws.Range("B15").Value = "MACEDONIA"
$ws.Range("D15").Value = 44360.75
$ws.Range("D16").Value = 44364.875
$ws.Range("B17").Value = "MACEDONIA"
$ws.Range("D17").Value = 44364.625
$ws.Range("A18").Value = "MACEDONIA"
$ws.Range("D18").Value = 44368.75
$ws.Range("D19").Value = 44368.75

$ws.Range("D20").Value = 44360.625
$ws.Range("A21").Value = "SCOTLAND"
$ws.Range("D21").Value = 44361.625
$ws.Range("B22").Value = "SCOTLAND"
$ws.Range("D22").Value = 44365.875
$ws.Range("D23").Value = 44365.75
$ws.Range("D24").Value = 44369.875
$ws.Range("B25").Value = "SCOTLAND"
$ws.Range("D25").Value = 44369.875

$ws.Range("D26").Value = 44361.875
$ws.Range("B27").Value = "SLOVAKIA"
$ws.Range("D27").Value = 44361.75
$ws.Range("B28").Value = "SLOVAKIA"
$ws.Range("D28").Value = 44365.625
$ws.Range("D29").Value = 44366.875
$ws.Range("A30").Value = "SLOVAKIA"
$ws.Range("D30").Value = 44370.75
$ws.Range("D31").Value = 44370.75

$ws.Range("D32").Value = 44362.875
$ws.Range("A33").Value = "HUNGARY"
$ws.Range("D33").Value = 44362.75
$ws.Range("D34").Value = 44366.75
$ws.Range("A35").Value = "HUNGARY"
$ws.Range("D35").Value = 44366.625
$ws.Range("B36").Value = "HUNGARY"
$ws.Range("D36").Value = 44370.875
$ws.Range("D37").Value = 44370.875

# Leave the Matches sheet active/selected, matching the saved view.
$ws.Activate()
$ws.Range("B16").Select()
